# Apply the edits described by the commit "Image width and height".
$p = $ppt.ActivePresentation

# --- Slide 1: update title / subtitle text -------------------------------
$s1 = $p.Slides.Item(1)
$titleTr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$titleTr1.Characters(1, $titleTr1.Length).Text = "Hello World! example"

$subTr1 = $s1.Shapes.Item(2).TextFrame.TextRange
$subTr1.InsertAfter("python-pptx was here!") | Out-Null

# --- Slide 2: update title + bullet paragraphs with outline levels -------
$s2 = $p.Slides.Item(2)
$titleTr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$titleTr2.Characters(1, $titleTr2.Length).Text = "Adding a Bullet Slide"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2).Text = "Use h3 for bullet slide layout"
$tr2.Paragraphs(3).Text = "Use h4 for first bullet"
$tr2.Paragraphs(3).IndentLevel = 2
$tr2.Paragraphs(4).Text = "Use h5 for subsequent bullet"
$tr2.Paragraphs(4).IndentLevel = 3

# --- Slide 3: remove all shapes (title + content placeholder) ------------
$s3 = $p.Slides.Item(3)
$guard = 0
while ($s3.Shapes.Count -gt 0 -and $guard -lt 10) {
    $s3.Shapes.Item(1).Delete()
    $guard = $guard + 1
}

# --- Old slide 5 (picture slide): shrink first picture, add a second -----
$s5 = $p.Slides.Item(5)
$pic1 = $s5.Shapes.Item(1)
$pic1.Width = 72
$pic1.Height = 72

$dup = $pic1.Duplicate()
$pic2 = $dup.Item(1)
$pic2.Name = "Picture 2"
$pic2.Left = 216
$pic2.Top = 72
$pic2.Width = 400
$pic2.Height = 400

# --- Remove old slide 4 (TextBox slide) so former slide 5 becomes 4 ------
$p.Slides.Item(4).Delete()
